$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("V5").Value = 13
$ws.Range("V6").Value = 4
$ws.Range("U7").Value = 8
$ws.Range("V16").Value = 18
$ws.Range("U20").Value = 15
$ws.Range("U27").Value = 21
$ws.Range("U28").Value = 0
$ws.Range("V28").Value = 29
$ws.Range("U32").Value = 2
$ws.Range("U35").Value = 10
$ws.Range("V42").Value = 8
$ws.Range("U44").Value = 1
$ws.Range("V44").Value = 23
$ws.Range("U48").Value = 2
$ws.Range("V48").Value = 25
$ws.Range("U56").Value = 13
$ws.Range("U57").Value = 2
$ws.Range("U58").Value = 19
$ws.Range("V58").Value = 19
$ws.Range("U60").Value = 30
$ws.Range("U61").Value = 28
$ws.Range("U63").Value = 0
$ws.Range("V63").Value = 54
$ws.Range("U68").Value = 7
$ws.Range("V68").Value = 27
$ws.Range("U70").Value = 1
$ws.Range("V70").Value = 42
$ws.Range("U71").Value = 21
$ws.Range("U73").Value = 8
$ws.Range("V75").Value = 12
$ws.Range("U76").Value = 1
$ws.Range("V76").Value = 11
$ws.Range("U77").Value = 20
$ws.Range("V77").Value = 12
$ws.Range("U82").Value = 35
$ws.Range("V82").Value = 15
$ws.Range("U83").Value = 20
$ws.Range("U84").Value = 27
$ws.Range("V84").Value = 18
$ws.Range("U85").Value = 0
$ws.Range("V85").Value = 35
$ws.Range("V86").Value = 19
$ws.Range("U87").Value = 20
$ws.Range("U88").Value = 17
$ws.Range("U92").Value = 61
$ws.Range("U98").Value = 17
$ws.Range("V98").Value = 15
$ws.Range("U101").Value = 37
$ws.Range("U106").Value = 15
$ws.Range("V115").Value = 23
$ws.Range("U118").Value = 20
$ws.Range("U120").Value = 32
$ws.Range("U123").Value = 34
$ws.Range("U130").Value = 36
$ws.Range("U142").Value = 0
$ws.Range("U145").Value = 0
$ws.Range("V145").Value = 25
$ws.Range("U152").Value = 9
$ws.Range("V152").Value = 21
$ws.Range("U153").Value = 13
$ws.Range("U156").Value = 67
$ws.Range("U158").Value = 0
$ws.Range("V158").Value = 1
$ws.Range("U159").Value = 11
$ws.Range("U162").Value = 0
$ws.Range("V162").Value = 11
$ws.Range("U163").Value = 0
$ws.Range("U167").Value = 4
$ws.Range("U168").Value = 3
$ws.Range("V168").Value = 6
$ws.Range("U175").Value = 11
$ws.Range("U177").Value = 8
$ws.Range("U182").Value = 4
$ws.Range("V182").Value = 42
$ws.Range("U184").Value = 33
$ws.Range("U185").Value = 29
$ws.Range("U186").Value = 3
$ws.Range("V186").Value = 12
$ws.Range("U187").Value = 7
$ws.Range("U193").Value = 2
$ws.Range("V193").Value = 21
$ws.Range("U194").Value = 8
$ws.Range("U197").Value = 0
$ws.Range("V197").Value = 47
$ws.Range("U203").Value = 14
$ws.Range("U206").Value = 5
$ws.Range("V206").Value = 18
$ws.Range("U209").Value = 33
$ws.Range("V214").Value = 12
$ws.Range("U222").Value = 9
$ws.Range("U223").Value = 30
$ws.Range("U228").Value = 0
$ws.Range("V228").Value = 8
$ws.Range("U234").Value = 22
$ws.Range("U237").Value = 8
$ws.Range("U252").Value = 6
$ws.Range("V253").Value = 30
